# issue #5: stock data from json to db
#
# The 股票 (stock) sheet gains three new columns: "category", "source_file"
# and "index". "category" is inserted right after "property_category" (so
# the existing "date", "legislator_name" and "legislator_id" columns shift
# one column to the right), while "source_file" and "index" are appended
# after "legislator_id".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票

$xlPasteValues  = -4163
$xlPasteFormats = -4122

# ------------------------------------------------------------------
# Row 1 (header row): stamp the existing header style (K1) onto the
# three brand-new header cells first, so the new cells end up sharing
# the same style id as the rest of the header row, then shift the
# existing header labels one column to the right and add the new ones.
# ------------------------------------------------------------------
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial($xlPasteFormats)
$ws.Range("M1").PasteSpecial($xlPasteFormats)
$ws.Range("N1").PasteSpecial($xlPasteFormats)

$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial($xlPasteValues)   # legislator_id (was K1)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteValues)   # legislator_name (was J1)
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial($xlPasteValues)   # date (was I1)

$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ------------------------------------------------------------------
# Row 2 (data row): same idea - copy the existing data-cell format onto
# the new cells, then shift the data values right and fill in the new
# ones.
# ------------------------------------------------------------------
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial($xlPasteFormats)
$ws.Range("M2").PasteSpecial($xlPasteFormats)
$ws.Range("N2").PasteSpecial($xlPasteFormats)

$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial($xlPasteValues)   # legislator_id value (was K2, 919)
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial($xlPasteValues)   # legislator_name value (was J2)
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial($xlPasteValues)   # date value (was I2)

$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmpc261"
$ws.Range("N2").Value = 68
